# Update countries & provincias Spain
# Applies the data refresh captured in the target diff:
#  - Timestamp in A1 bumped from 05:42 to 06:59
#  - Three countries (Kirguistan / Republica de Macedonia / Venezuela) re-sorted,
#    each row keeping its position but getting the next country's data (cascading update)
#    and Kirguistan (now row 81) receiving freshly updated case counts.
#  - Three pairs of tied-count countries swap display order (Laos/Santa Lucia,
#    Dominica/Fiyi, Islas Malvinas/Groenlandia) - numbers unchanged, only labels swap.
#  - Pakistan (row 15) and Tailandia (row 99) get updated case counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 06:59"

# --- Kenia / Macedonia / Venezuela / Kirguistan / Haiti block (rows 80-84) ---
# Row 81: Macedonia -> Kirguistan, with fresh data
$ws.Range("A81").Value = "Kirguistan"
$ws.Range("B81").Value = 6767
$ws.Range("C81").Value = 506
$ws.Range("D81").Value = 2655
$ws.Range("E81").Value = 4036
$ws.Range("G81").Value = 10
$ws.Range("H81").Value = 76

# Row 82: Venezuela -> Macedonia, data shifted from old row 81
$ws.Range("A82").Value = "Republica de Macedonia"
$ws.Range("B82").Value = 6625
$ws.Range("D82").Value = 2748
$ws.Range("E82").Value = 3556
$ws.Range("H82").Value = 321

# Row 83: Kirguistan -> Venezuela, data shifted from old row 82
$ws.Range("A83").Value = "Venezuela"
$ws.Range("B83").Value = 6273
$ws.Range("D83").Value = 2100
$ws.Range("E83").Value = 4116
$ws.Range("H83").Value = 57

# --- Nueva Caledonia / Laos / Santa Lucia / Dominica / Fiyi block (rows 202-211) ---
# Tied case counts, only display order (label) swaps; underlying numbers untouched.
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A204").Value = "Laos"

$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"

$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Pakistan (row 15) updated case counts ---
$ws.Range("B15").Value = 221896
$ws.Range("C15").Value = 4087
$ws.Range("D15").Value = 113623
$ws.Range("E15").Value = 103722
$ws.Range("G15").Value = 78
$ws.Range("H15").Value = 4551

# --- Tailandia (row 99) updated case counts ---
$ws.Range("B99").Value = 3180
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 3066
$ws.Range("E99").Value = 56
